$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 106
$ws.Range("H106").Value = 74078780
$ws.Range("I106").Value = 20838004
$ws.Range("K106").Value = 20838004
$ws.Range("M106").Value = -20837373
# Row 113
$ws.Range("H113").Value = 4617.8945
$ws.Range("I113").Value = 3738.25
$ws.Range("J113").Value = 6125.857
$ws.Range("K113").Value = 3738.25
$ws.Range("L113").Value = 6125.857
$ws.Range("M113").Value = -484.25
$ws.Range("N113").Value = -12633.857
# Row 137
$ws.Range("H137").Value = 1635.2667
$ws.Range("I137").Value = 1131.6923
$ws.Range("J137").Value = 2020.3529
$ws.Range("K137").Value = 3395.0769
$ws.Range("L137").Value = 6061.0587
$ws.Range("M137").Value = -845.0769
$ws.Range("N137").Value = -11161.0587
# Row 138
$ws.Range("H138").Value = 4504.4307
$ws.Range("I138").Value = 1379.2424
$ws.Range("J138").Value = 7148.8203
$ws.Range("K138").Value = 4137.7272
$ws.Range("L138").Value = 21446.4609
$ws.Range("M138").Value = 1002.2728
$ws.Range("N138").Value = -31726.4609

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1131.6086
$ws.Range("I2").Value = 1231.3125
$ws.Range("J2").Value = 903.7143
$ws.Range("K2").Value = 1231.3125
$ws.Range("L2").Value = 903.7143
$ws.Range("M2").Value = -1118.3125
$ws.Range("N2").Value = -1129.7143
# Row 32
$ws.Range("H32").Value = 5304.46
$ws.Range("I32").Value = 4816.0264
$ws.Range("J32").Value = 6851.1665
$ws.Range("K32").Value = 4816.0264
$ws.Range("L32").Value = 6851.1665
$ws.Range("M32").Value = -4529.0264
$ws.Range("N32").Value = -7425.1665
# Row 61
$ws.Range("H61").Value = 1659.5178
$ws.Range("I61").Value = 1525.625
$ws.Range("K61").Value = 1525.625
$ws.Range("M61").Value = -1313.625
# Row 74
$ws.Range("H74").Value = 1354.4082
$ws.Range("I74").Value = 1103.625
$ws.Range("J74").Value = 1826.4706
$ws.Range("K74").Value = 1103.625
$ws.Range("L74").Value = 1826.4706
$ws.Range("M74").Value = -229.625
$ws.Range("N74").Value = -3574.4706
# Row 77
$ws.Range("H77").Value = 1354.4082
$ws.Range("I77").Value = 1103.625
$ws.Range("J77").Value = 1826.4706
$ws.Range("K77").Value = 5518.125
$ws.Range("L77").Value = 9132.353000000001
$ws.Range("M77").Value = -1150.125
$ws.Range("N77").Value = -17868.353
# Row 109
$ws.Range("H109").Value = 35277
$ws.Range("J109").Value = 35277
$ws.Range("L109").Value = 35277
$ws.Range("N109").Value = -38051
# Row 116
$ws.Range("H116").Value = 1131.6086
$ws.Range("I116").Value = 1231.3125
$ws.Range("J116").Value = 903.7143
$ws.Range("K116").Value = 1231.3125
$ws.Range("L116").Value = 903.7143
$ws.Range("M116").Value = 1062.6875
$ws.Range("N116").Value = -5491.7143
# Row 123
$ws.Range("H123").Value = 31000
$ws.Range("J123").Value = 31000
$ws.Range("L123").Value = 31000
$ws.Range("N123").Value = -40800
# Row 125
$ws.Range("H125").Value = 47600
$ws.Range("J125").Value = 47600
$ws.Range("L125").Value = 47600
$ws.Range("N125").Value = -57440
# Row 132
$ws.Range("H132").Value = 2192.4324
$ws.Range("I132").Value = 1650.6471
$ws.Range("J132").Value = 8332.666999999999
$ws.Range("K132").Value = 4951.9413
$ws.Range("L132").Value = 24998.001
$ws.Range("M132").Value = -2421.9413
$ws.Range("N132").Value = -30058.001
# Row 136
$ws.Range("H136").Value = 1659.5178
$ws.Range("I136").Value = 1525.625
$ws.Range("K136").Value = 4576.875
$ws.Range("M136").Value = -2026.875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1131.6086
$ws.Range("I3").Value = 1231.3125
$ws.Range("J3").Value = 903.7143
$ws.Range("K3").Value = 1231.3125
$ws.Range("L3").Value = 903.7143
$ws.Range("M3").Value = -1117.3125
$ws.Range("N3").Value = -1131.7143
# Row 107
$ws.Range("H107").Value = 1399.25
$ws.Range("I107").Value = 1260.625
$ws.Range("J107").Value = 1676.5
$ws.Range("K107").Value = 1260.625
$ws.Range("L107").Value = 1676.5
$ws.Range("M107").Value = 659.375
$ws.Range("N107").Value = -5516.5
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# Row 134
$ws.Range("H134").Value = 2113.8108
$ws.Range("I134").Value = 1375.909
$ws.Range("K134").Value = 4127.727000000001
$ws.Range("M134").Value = -1592.727000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4034.9556
$ws.Range("I31").Value = 4208.467
$ws.Range("J31").Value = 3948.2
$ws.Range("K31").Value = 4208.467
$ws.Range("L31").Value = 3948.2
$ws.Range("M31").Value = -3913.467
$ws.Range("N31").Value = -4538.2
# Row 34
$ws.Range("H34").Value = 4034.9556
$ws.Range("I34").Value = 4208.467
$ws.Range("J34").Value = 3948.2
$ws.Range("K34").Value = 4208.467
$ws.Range("L34").Value = 3948.2
$ws.Range("M34").Value = -4006.467
$ws.Range("N34").Value = -4352.2
# Row 62
$ws.Range("H62").Value = 37041668
$ws.Range("I62").Value = 2725
$ws.Range("J62").Value = 66672820
$ws.Range("K62").Value = 2725
$ws.Range("L62").Value = 66672820
$ws.Range("M62").Value = -2101
$ws.Range("N62").Value = -66674068
# Row 65
$ws.Range("H65").Value = 37041668
$ws.Range("I65").Value = 2725
$ws.Range("J65").Value = 66672820
$ws.Range("K65").Value = 13625
$ws.Range("L65").Value = 333364100
$ws.Range("M65").Value = -10505
$ws.Range("N65").Value = -333370340
# Row 122
$ws.Range("H122").Value = 1878.3
$ws.Range("I122").Value = 1138
$ws.Range("J122").Value = 2195.5715
$ws.Range("K122").Value = 3414
$ws.Range("L122").Value = 6586.7145
$ws.Range("M122").Value = -964
$ws.Range("N122").Value = -11486.7145
# Row 132
$ws.Range("H132").Value = 3422.394
$ws.Range("I132").Value = 3068.8215
$ws.Range("J132").Value = 5402.4
$ws.Range("K132").Value = 9206.4645
$ws.Range("L132").Value = 16207.2
$ws.Range("M132").Value = -6676.4645
$ws.Range("N132").Value = -21267.2
# Row 134
$ws.Range("H134").Value = 4623.391
$ws.Range("I134").Value = 5494.4
$ws.Range("J134").Value = 3953.3845
$ws.Range("K134").Value = 16483.2
$ws.Range("L134").Value = 11860.1535
$ws.Range("M134").Value = -13948.2
$ws.Range("N134").Value = -16930.1535

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 3873.8462
$ws.Range("I68").Value = 6660.222
$ws.Range("J68").Value = 1485.5238
$ws.Range("K68").Value = 19980.666
$ws.Range("L68").Value = 4456.5714
$ws.Range("M68").Value = -19169.666
$ws.Range("N68").Value = -6078.5714
# Row 71
$ws.Range("H71").Value = 3873.8462
$ws.Range("I71").Value = 6660.222
$ws.Range("J71").Value = 1485.5238
$ws.Range("K71").Value = 59941.998
$ws.Range("L71").Value = 13369.7142
$ws.Range("M71").Value = -55885.998
$ws.Range("N71").Value = -21481.7142
# Row 113
$ws.Range("H113").Value = 1538939.5
$ws.Range("I113").Value = 2174362.8
$ws.Range("J113").Value = 625518.4
$ws.Range("K113").Value = 6523088.399999999
$ws.Range("L113").Value = 1876555.2
$ws.Range("M113").Value = -6520918.399999999
$ws.Range("N113").Value = -1880895.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1411.6364
$ws.Range("I102").Value = 1125
$ws.Range("J102").Value = 2176
$ws.Range("K102").Value = 1125
$ws.Range("L102").Value = 2176
$ws.Range("M102").Value = 497
$ws.Range("N102").Value = -5420
# Row 126
$ws.Range("H126").Value = 8071.032
$ws.Range("I126").Value = 9842.166999999999
$ws.Range("K126").Value = 29526.501
$ws.Range("M126").Value = -27056.501
# Row 132
$ws.Range("H132").Value = 2903.1
$ws.Range("I132").Value = 2813.12
$ws.Range("K132").Value = 8439.360000000001
$ws.Range("M132").Value = -5909.360000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 17857524
$ws.Range("I55").Value = 259.23077
$ws.Range("J55").Value = 33333820
$ws.Range("K55").Value = 259.23077
$ws.Range("L55").Value = 33333820
$ws.Range("M55").Value = -86.23077000000001
$ws.Range("N55").Value = -33334166
# Row 132
$ws.Range("H132").Value = 11114364
$ws.Range("I132").Value = 14945668
$ws.Range("J132").Value = 3584.8
$ws.Range("K132").Value = 44837004
$ws.Range("L132").Value = 10754.4
$ws.Range("M132").Value = -44834474
$ws.Range("N132").Value = -15814.4
# Row 136
$ws.Range("H136").Value = 5082.737
$ws.Range("I136").Value = 2563.0356
$ws.Range("J136").Value = 12137.9
$ws.Range("K136").Value = 7689.1068
$ws.Range("L136").Value = 36413.7
$ws.Range("M136").Value = -5139.1068
$ws.Range("N136").Value = -41513.7

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 108
$ws.Range("H108").Value = 40454.547
$ws.Range("J108").Value = 40454.547
$ws.Range("L108").Value = 40454.547
$ws.Range("N108").Value = -48134.547
# Row 132
$ws.Range("H132").Value = 2197.0356
$ws.Range("I132").Value = 2032.6818
$ws.Range("J132").Value = 2799.6667
$ws.Range("K132").Value = 6098.0454
$ws.Range("L132").Value = 8399.000100000001
$ws.Range("M132").Value = -3568.0454
$ws.Range("N132").Value = -13459.0001
# Row 136
$ws.Range("H136").Value = 932.4878
$ws.Range("I136").Value = 772
$ws.Range("J136").Value = 1430
$ws.Range("K136").Value = 2316
$ws.Range("L136").Value = 4290
$ws.Range("M136").Value = 234
$ws.Range("N136").Value = -9390
